# sample_employeeBulkQuickOnboarding.xlsx - bulk/quick onboarding template tweaks
# - rename the "Reporting Manager Code" header to "L1 Manager Code"
# - refresh the sample email (value + hyperlink) and sample DOJ-format placeholder
# - refresh the sample mobile number

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DOJ sample text: "YY/MM/DD" -> "dd-mm-yyyy"
$ws.Range("D2").Value = "dd-mm-yyyy"

# Email sample: update cell text and the underlying mailto hyperlink together
# (drop the old hyperlink, write the new text, re-add the hyperlink) so the
# link target and the displayed text stay in sync, then restore the cell's
# Hyperlink style.
$ws.Hyperlinks.Delete()
$ws.Range("C2").Value = "test@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:test@gmail.com")
$ws.Range("C2").Style = "Hyperlink"

# Header rename: "Reporting Manager Code " -> "L1 Manager Code"
$ws.Range("G1").Value = "L1 Manager Code"

# Mobile number sample value
$ws.Range("E2").Value = 1234567890
